$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rows 58 and 59 had their match data swapped (Sudtirol-Modena and
#    Spezia-Brescia entries traded places). Columns A-E are identical
#    between the two rows already, so only F:V need to be exchanged.
# ------------------------------------------------------------------
$swapCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($col in $swapCols) {
    $addr58 = $col + "58"
    $addr59 = $col + "59"
    $v58 = $ws.Range($addr58).Value2
    $v59 = $ws.Range($addr59).Value2
    $ws.Range($addr58).Value = $v59
    $ws.Range($addr59).Value = $v58
}

# ------------------------------------------------------------------
# 2) Two new match rows were appended at the bottom (rows 117 & 118).
#    Clone the formatting of the last existing row (116) for the
#    styled columns (A = bold/bordered index, E = date-time number
#    format) before writing the new values.
# ------------------------------------------------------------------
$ws.Range("A116").Copy()
$ws.Range("A117").PasteSpecial(-4122)
$ws.Range("A118").PasteSpecial(-4122)

$ws.Range("E116").Copy()
$ws.Range("E117").PasteSpecial(-4122)
$ws.Range("E118").PasteSpecial(-4122)

# Row 117: Lecco 0 - 0 Spezia
$ws.Range("A117").Value = 116
$ws.Range("B117").Value = "italy"
$ws.Range("C117").Value = "serie-b"
$ws.Range("D117").Value = "2023-2024"
$ws.Range("E117").Value = 45238.77083333334
$ws.Range("F117").Value = "Lecco"
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = "Spezia"
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 3.39
$ws.Range("K117").Value = "06/11/2023 15:12"
$ws.Range("L117").Value = 3.98
$ws.Range("M117").Value = "08/11/2023 18:26"
$ws.Range("N117").Value = 3.39
$ws.Range("O117").Value = "06/11/2023 15:12"
$ws.Range("P117").Value = 3.52
$ws.Range("Q117").Value = "08/11/2023 18:20"
$ws.Range("R117").Value = 2.26
$ws.Range("S117").Value = "06/11/2023 15:12"
$ws.Range("T117").Value = 2.03
$ws.Range("U117").Value = "08/11/2023 18:26"
$ws.Range("V117").Value = "https://www.betexplorer.com/football/italy/serie-b/lecco-spezia/GhYhP1tg/"

# Row 118: Palermo 1 - 0 Brescia
$ws.Range("A118").Value = 117
$ws.Range("B118").Value = "italy"
$ws.Range("C118").Value = "serie-b"
$ws.Range("D118").Value = "2023-2024"
$ws.Range("E118").Value = 45238.77083333334
$ws.Range("F118").Value = "Palermo"
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = "Brescia"
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 1.71
$ws.Range("K118").Value = "04/11/2023 17:13"
$ws.Range("L118").Value = 1.66
$ws.Range("M118").Value = "08/11/2023 18:24"
$ws.Range("N118").Value = 3.9
$ws.Range("O118").Value = "04/11/2023 17:13"
$ws.Range("P118").Value = 3.68
$ws.Range("Q118").Value = "08/11/2023 18:24"
$ws.Range("R118").Value = 5.12
$ws.Range("S118").Value = "04/11/2023 17:13"
$ws.Range("T118").Value = 6.36
$ws.Range("U118").Value = "08/11/2023 18:24"
$ws.Range("V118").Value = "https://www.betexplorer.com/football/italy/serie-b/palermo-brescia/Q5CsOLOl/"
